$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Area 5 item sprites (columns C/D) for rows 9-13, mirroring the
#     existing "item type / item letter" pairs already present in rows 4-8 ---
$ws.Range("C9").Value  = "pinetree"
$ws.Range("D9").Value  = "F"

$ws.Range("C10").Value = "snowrock"
$ws.Range("D10").Value = "Q"

$ws.Range("C11").Value = "snow"
$ws.Range("D11").Value = "N"

$ws.Range("C12").Value = "mudwall"
$ws.Range("D12").Value = "V"

$ws.Range("C13").Value = "ice"
$ws.Range("D13").Value = "I"

# --- Area 5 tile sprites (column I) pick up the bold formatting that the
#     rest of the "blocks" letter column already uses ---
$ws.Range("I9").Font.Bold  = $true
$ws.Range("I12").Font.Bold = $true
$ws.Range("I17").Font.Bold = $true
$ws.Range("I20").Font.Bold = $true
$ws.Range("I25").Font.Bold = $true

# --- restore the selection to where the author left off editing ---
$ws.Range("K20").Select() | Out-Null
